# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos list price/volume update described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price text (must stay plain text,
# e.g. '0.9980' must not collapse to 0.998) -> force text format first.
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '29.324.18'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.844.16'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").Value = '0.9976'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '240.12'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").Value = '0.6273'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '0.9986'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.07489'
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").Value = '0.2899'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '24.41'
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("D11").Value = '0.07734'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '1.845.82'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '4.985'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '0.6795'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '0.00001050'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '82.03'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").Value = '6.179'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '29.382.52'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = '229.37'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '0.9994'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '7.487'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").Value = '0.9993'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '158.58'
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").Value = '8.415'
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '17.52'
$ws.Range("D28").Value = '0.06400'
$ws.Range("E28").Value = '  +14.33%  '
$ws.Range("D29").Value = '1.424'
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("D30").Value = '1.482'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").Value = '4.088'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").Value = '4.094'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = '1.833'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '1.141'
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").Value = '0.6976'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '2.578'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '1.269.34'
$ws.Range("E37").Value = '  +3.55%  '
$ws.Range("D38").Value = '2.835'
$ws.Range("E38").Value = '  +4.29%  '
$ws.Range("D39").Value = '0.01832'
$ws.Range("E39").Value = '  +1.77%  '
$ws.Range("D40").Value = '6.709'
$ws.Range("E40").Value = '  +5.63%  '
$ws.Range("D41").Value = '0.9139'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").Value = '0.9980'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '2.009.50'
$ws.Range("E43").Value = '  -18.34%  '
$ws.Range("D44").Value = '101.18'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '66.23'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D46").Value = '1.731'
$ws.Range("E46").Value = '  +3.45%  '
$ws.Range("D47").Value = '7.072'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").Value = '0.1164'
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("D49").Value = '9.016'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").Value = '0.3961'
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").Value = '0.00000000114'
$ws.Range("E51").Value = '  -3.58%  '
